# Weekly fruit/vegetable price update: 4 new weekly records were added for
# "Repollo" (Vega Central Mapocho de Santiago) at the top of the existing
# history block (rows 727-742), pushing the prior rows down by 4 (now
# 731-746). The sheet's used range grows from A1:R742 to A1:R746.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 727; Excel shifts rows 727:742 down to
# 731:746 and copies formatting (e.g. the date style on column D) from the
# row above, same as a normal in-app row insert.
$ws.Rows("727:730").Insert()

# Fill in the 4 new rows with this week's data.
$newRows = @(
    @{ Row = 727; H = "Crespo record"; I = "Primera"; J = 3400; K = 1200; L = 1300; M = 1250; O = "Provincia de Quillota"; P = 1250 },
    @{ Row = 728; H = "Crespo record"; I = "Segunda"; J = 1690; K = 1000; L = 1000; M = 1000; O = "Provincia de Quillota"; P = 1000 },
    @{ Row = 729; H = "Morada(o)";     I = "Primera"; J = 1600; K = 1300; L = 1500; M = 1400; O = "Provincia de Quillota"; P = 1400 },
    @{ Row = 730; H = "Morada(o)";     I = "Segunda"; J = 790;  K = 1100; L = 1100; M = 1100; O = "Provincia de Quillota"; P = 1100 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 45041
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112006
    $ws.Cells.Item($row, 7).Value = "Repollo"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
